# Fault Calcs.xlsx - "Worked more on fault analysis"
# Applies the update described by the commit diff:
#  - Sheet "Part A": recompute/replace the post-fault currents for buses 4-17
#    (F7:G20) with new values, converts the D4:D20 magnitude formulas into a
#    shared-formula group, and moves the active selection.
#  - Sheet "Part B": fills in the still-missing Real/Imag/Mag (I0+/I-) columns
#    for "Bus 2" and "Bus 4" (B:H), refreshes the J7:K7 fault currents to match
#    the new Part A numbers, and appends a brand-new small summary table
#    (rows 9-12) with headers I0 / I+ / I- and per-bus Real/Angle values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Part A"
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Part A")

# New fault-current (F/G) values for bus rows 7-20; column H recomputes via
# the existing shared SQRT formula automatically.
$wsA.Range("F7").Value = 3.9594999999999998
$wsA.Range("G7").Value = -13.8637

$wsA.Range("F8").Value = 6.6231999999999998
$wsA.Range("G8").Value = -20.9344

$wsA.Range("F9").Value = 4.9508000000000001
$wsA.Range("G9").Value = -14.007199999999999

$wsA.Range("F10").Value = 5.0475000000000003
$wsA.Range("G10").Value = -15.167

$wsA.Range("F11").Value = 5.4569999999999999
$wsA.Range("G11").Value = -15.664

$wsA.Range("F12").Value = 2.5451999999999999
$wsA.Range("G12").Value = -12.442

$wsA.Range("F13").Value = 0.58919999999999995
$wsA.Range("G13").Value = -2.343

$wsA.Range("F14").Value = 5.1891999999999996
$wsA.Range("G14").Value = -18.247900000000001

$wsA.Range("F15").Value = 6.9977999999999998
$wsA.Range("G15").Value = -18.9666

$wsA.Range("F16").Value = 0.60450000000000004
$wsA.Range("G16").Value = -2.1263999999999998

$wsA.Range("F17").Value = 3.5808
$wsA.Range("G17").Value = -13.5433

$wsA.Range("F18").Value = 4.5118999999999998
$wsA.Range("G18").Value = -12.6579

$wsA.Range("F19").Value = 1.2270000000000001
$wsA.Range("G19").Value = -6.3997000000000002

$wsA.Range("F20").Value = 0.53620000000000001
$wsA.Range("G20").Value = -6.4583000000000004

# Re-enter the magnitude formulas for D4:D20 so they become one shared group.
$wsA.Range("D4:D20").Formula = "=SQRT((B4)^2+(C4)^2)"

# Move the active selection on this sheet.
[void]$wsA.Activate()
[void]$wsA.Range("C45").Select()

# ---------------------------------------------------------------------------
# Sheet "Part B"
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Part B")

# Fill in the previously-empty Real/Imag/Mag block (columns B:D and F:H) for
# the "Bus 2" row (row 6) and the "Bus 4" row (row 7).
$wsB.Range("B6").Value = -46.525599999999997
$wsB.Range("C6").Value = -13.553699999999999
$wsB.Range("D6").Formula = "=SQRT((B6)^2+(C6)^2)"

$wsB.Range("F6").Value = -56.765700000000002
$wsB.Range("G6").Value = 25.816299999999998
$wsB.Range("H6").Formula = "=SQRT((F6)^2+(G6)^2)"

$wsB.Range("B7").Value = -18.116700000000002
$wsB.Range("C7").Value = -5.7504
$wsB.Range("D7").Formula = "=SQRT((B7)^2+(C7)^2)"

$wsB.Range("F7").Value = -19.6312
$wsB.Range("G7").Value = 0.14019999999999999
$wsB.Range("H7").Formula = "=SQRT((F7)^2+(G7)^2)"

# Updated SLG fault values (mirrors the new Part A bus-4 row).
$wsB.Range("J7").Value = 3.9594999999999998
$wsB.Range("K7").Value = -13.8637

# New summary table: headers.
$wsB.Range("B9").Value = "I0"
$wsB.Range("D9").Value = "I+"
$wsB.Range("F9").Value = "I-"

$wsB.Range("B10").Value = "Mag"
$wsB.Range("C10").Value = "Angle"
$wsB.Range("D10").Value = "Mag"
$wsB.Range("E10").Value = "Angle"
$wsB.Range("F10").Value = "Mag"
$wsB.Range("G10").Value = "Angle"

# New summary table: data rows.
$wsB.Range("A11").Value = "Bus 2"
$wsB.Range("B11").Value = 22.4
$wsB.Range("C11").Value = -74.435500000000005
$wsB.Range("D11").Value = 22.4
$wsB.Range("E11").Value = -74.435500000000005
$wsB.Range("F11").Value = 22.4
$wsB.Range("G11").Value = -74.435500000000005

$wsB.Range("A12").Value = "Bus 4"
$wsB.Range("B12").Value = 5.2195299999999998
$wsB.Range("C12").Value = -74.442599999999999
$wsB.Range("D12").Value = 5.2195299999999998
$wsB.Range("E12").Value = -74.442599999999999
$wsB.Range("F12").Value = 5.2195299999999998
$wsB.Range("G12").Value = -74.442599999999999

# Move the active selection/active sheet to match the saved workbook state.
[void]$wsB.Activate()
[void]$wsB.Range("A8").Select()
